$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" (column F) values for the rows that were re-pulled.
$ws.Range("F2").Value  = -7
$ws.Range("F3").Value  = -2
$ws.Range("F4").Value  = -1
$ws.Range("F5").Value  = 0
$ws.Range("F10").Value = -4
$ws.Range("F17").Value = -5
$ws.Range("F22").Value = 1
$ws.Range("F23").Value = -5
$ws.Range("F25").Value = -3
$ws.Range("F34").Value = 4
$ws.Range("F37").Value = -4
$ws.Range("F38").Value = -3
$ws.Range("F40").Value = -8
